$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44162
$ws.Range("K2").Value = 1000
$ws.Range("L2").Value = 1000
$ws.Range("M2").Value = 1000
$ws.Range("N2").Value = "`$/atado"
$ws.Range("P2").Value = 1000

# Row 3
$ws.Range("D3").Value = 44455
$ws.Range("J3").Value = 1500
$ws.Range("K3").Value = 2400
$ws.Range("L3").Value = 2400
$ws.Range("M3").Value = 2400
$ws.Range("O3").Value = "Región del Maule"
$ws.Range("P3").Value = 2400

# Row 4
$ws.Range("D4").Value = 44176
$ws.Range("J4").Value = 2000
$ws.Range("K4").Value = 900
$ws.Range("L4").Value = 900
$ws.Range("M4").Value = 900
$ws.Range("O4").Value = "Provincia de Linares"
$ws.Range("P4").Value = 900

# Row 5
$ws.Range("D5").Value = 44459
$ws.Range("J5").Value = 1500
$ws.Range("K5").Value = 2000
$ws.Range("L5").Value = 2000
$ws.Range("M5").Value = 2000
$ws.Range("O5").Value = "Provincia de Limarí"
$ws.Range("P5").Value = 2000

# Row 6
$ws.Range("D6").Value = 44462
$ws.Range("J6").Value = 2000
$ws.Range("K6").Value = 1800
$ws.Range("L6").Value = 2000
$ws.Range("M6").Value = 1900
$ws.Range("P6").Value = 1900

# Row 7
$ws.Range("D7").Value = 44160
$ws.Range("K7").Value = 800
$ws.Range("L7").Value = 800
$ws.Range("M7").Value = 800
$ws.Range("P7").Value = 800

# Row 8
$ws.Range("D8").Value = 44460
$ws.Range("J8").Value = 2000
$ws.Range("K8").Value = 2000
$ws.Range("L8").Value = 2000
$ws.Range("M8").Value = 2000
$ws.Range("P8").Value = 2000

# Row 10
$ws.Range("D10").Value = 44175
$ws.Range("J10").Value = 800
$ws.Range("K10").Value = 1000
$ws.Range("L10").Value = 1100
$ws.Range("M10").Value = 1050
$ws.Range("O10").Value = "Provincia de Linares"
$ws.Range("P10").Value = 1050

# Row 11
$ws.Range("D11").Value = 44463

# Row 12
$ws.Range("D12").Value = 44169
$ws.Range("J12").Value = 3000
$ws.Range("K12").Value = 1000
$ws.Range("L12").Value = 1000
$ws.Range("M12").Value = 1000
$ws.Range("P12").Value = 1000

# Row 13
$ws.Range("D13").Value = 44467
$ws.Range("K13").Value = 1800
$ws.Range("L13").Value = 1800
$ws.Range("M13").Value = 1800
$ws.Range("O13").Value = "Provincia de Linares"
$ws.Range("P13").Value = 1800

# Row 14
$ws.Range("D14").Value = 44473
$ws.Range("J14").Value = 4000
$ws.Range("K14").Value = 1200
$ws.Range("L14").Value = 1200
$ws.Range("M14").Value = 1200
$ws.Range("P14").Value = 1200

# Row 15
$ws.Range("D15").Value = 44469
$ws.Range("K15").Value = 1200
$ws.Range("L15").Value = 1200
$ws.Range("M15").Value = 1200
$ws.Range("P15").Value = 1200

# Row 16
$ws.Range("D16").Value = 44168
$ws.Range("J16").Value = 3000
$ws.Range("K16").Value = 1000
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = 1000
$ws.Range("P16").Value = 1000

# Row 17
$ws.Range("D17").Value = 44161
$ws.Range("J17").Value = 3000

# Row 18
$ws.Range("D18").Value = 44165
$ws.Range("J18").Value = 2000
$ws.Range("K18").Value = 1200
$ws.Range("L18").Value = 1200
$ws.Range("M18").Value = 1200
$ws.Range("P18").Value = 1200

# Row 19
$ws.Range("D19").Value = 44475
$ws.Range("K19").Value = 1000
$ws.Range("L19").Value = 1100
$ws.Range("M19").Value = 1040
$ws.Range("P19").Value = 1040

# Row 20
$ws.Range("D20").Value = 44474
$ws.Range("J20").Value = 5000
$ws.Range("K20").Value = 1200
$ws.Range("L20").Value = 1200
$ws.Range("M20").Value = 1200
$ws.Range("P20").Value = 1200

# Row 21
$ws.Range("D21").Value = 44166
$ws.Range("J21").Value = 1500
$ws.Range("L21").Value = 1000
$ws.Range("M21").Value = 1000
$ws.Range("P21").Value = 1000

# Row 22
$ws.Range("D22").Value = 44468
$ws.Range("J22").Value = 3000
$ws.Range("K22").Value = 1500
$ws.Range("L22").Value = 1500
$ws.Range("M22").Value = 1500
$ws.Range("P22").Value = 1500

# Row 23
$ws.Range("D23").Value = 44172
$ws.Range("J23").Value = 2000
$ws.Range("N23").Value = "`$/kilo"
$ws.Range("O23").Value = "Región del Maule"

# Row 24
$ws.Range("D24").Value = 44466
$ws.Range("J24").Value = 2500
$ws.Range("K24").Value = 1800
$ws.Range("L24").Value = 1800
$ws.Range("M24").Value = 1800
$ws.Range("P24").Value = 1800

# New row 25
$ws.Range("D25").NumberFormat = $ws.Range("D24").NumberFormat
$ws.Range("A25").Value = 5
$ws.Range("B25").Value = "Macroferia Regional de Talca"
$ws.Range("C25").Value = "Maule"
$ws.Range("D25").Value = 44461
$ws.Range("E25").Value = 7
$ws.Range("F25").Value = 300000000
$ws.Range("G25").Value = "Espárragos"
$ws.Range("H25").Value = "Verde"
$ws.Range("I25").Value = "Primera"
$ws.Range("J25").Value = 2500
$ws.Range("K25").Value = 2000
$ws.Range("L25").Value = 2000
$ws.Range("M25").Value = 2000
$ws.Range("N25").Value = "`$/kilo"
$ws.Range("O25").Value = "Provincia de Linares"
$ws.Range("P25").Value = 2000
$ws.Range("Q25").Value = 1
$ws.Range("R25").Value = "Hortaliza"
